$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

$newUrl = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/communication-on-behalf-of"

# "Metadata" sheet: Property/Value table (rebrand Alvearie -> LinuxForHealth,
# bump version, refresh publish date).
$meta.Range("B2").Value = $newUrl               # URL
$meta.Range("B3").Value = "8.0.0"                # Version
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"  # Date
$meta.Range("B9").Value = "LinuxForHealth Team"  # Publisher

# The canonical URL is also referenced from the "Elements" sheet, as the
# fixed value of the Extension.url element (row 5, "Fixed Value" column) -
# keep it in sync with the Metadata sheet.
$elements.Range("Q5").Value = $newUrl
